$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.300.90"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "3.442.34"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'412.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "'130.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.69%  "
$ws.Range("E7").Value = "  +6.57%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "'0.756"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.90%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.20%  "
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("E13").Value = "  +6.68%  "
$ws.Range("D14").Value = "'20.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.94%  "
$ws.Range("D15").Value = "'0.0000196"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +53.73%  "
$ws.Range("D16").Value = "3.434.44"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("D17").Value = "'12.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +12.65%  "
$ws.Range("D18").Value = "'1.06"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.91%  "
$ws.Range("D19").Value = "62.238.78"
$ws.Range("E19").Value = "  +1.81%  "
$ws.Range("D20").Value = "'406.21"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +29.19%  "
$ws.Range("D21").Value = "'90.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.70%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("D23").Value = "'13.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("D24").Value = "'3.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.64%  "
$ws.Range("D25").Value = "'32.95"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.69%  "
$ws.Range("D26").Value = "'4.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").Value = "'8.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E28").Value = "  +3.11%  "
$ws.Range("E29").Value = "  +10.17%  "
$ws.Range("D30").Value = "'44.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.80%  "
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").Value = "'11.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +2.70%  "
$ws.Range("D36").Value = "'52.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("D38").Value = "'3.41"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").Value = "'2.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E40").Value = "  +7.18%  "
$ws.Range("D41").Value = "'0.315"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.15%  "
$ws.Range("D42").Value = "'140.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.97%  "
$ws.Range("D43").Value = "'1.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.21%  "
$ws.Range("E44").Value = "  +0.54%  "
$ws.Range("D45").Value = "'16.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.14%  "
$ws.Range("E46").Value = "  +4.61%  "
$ws.Range("D47").Value = "'22.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.52%  "
$ws.Range("D48").Value = "2.125.02"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +8.26%  "
